$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Row 36: add date, activity text (existing shared string) and hours
$ws.Range("A36").Value = 41443
$ws.Range("B36").Value = "LegacyIDE.e4xmi, Fragment und Prozessor Ansatz ausprobiert, Versionprobleme"
$ws.Range("C36").Value = 3

# Row 37: add date and new activity text (new shared string)
$ws.Range("A37").Value = 41444
$ws.Range("B37").Value = "Projekthandbuch, Statusbericht, Projektbericht"

# Update the current selection to match the authored state
$ws.Activate()
$ws.Range("B37").Select()
